$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.405.89"
$ws.Range("E2").Value = "  +12.53%  "

# Row 3
$ws.Range("D3").Value = "1.817.92"
$ws.Range("E3").Value = "  +7.66%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.68"
$ws.Range("E5").Value = "  +4.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.547"
$ws.Range("E6").Value = "  +4.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.49"
$ws.Range("E8").Value = "  +2.13%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.78"
$ws.Range("E9").Value = "  +3.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.283"
$ws.Range("E10").Value = "  +6.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0681"
$ws.Range("E11").Value = "  +8.96%  "

# Row 13
$ws.Range("D13").Value = "2.081.41"
$ws.Range("E13").Value = "  +7.74%  "

# Row 14
$ws.Range("D14").Value = "1.835.94"
$ws.Range("E14").Value = "  +8.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +3.94%  "

# Row 16
$ws.Range("D16").Value = "34.340.29"
$ws.Range("E16").Value = "  +12.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.24"
$ws.Range("E17").Value = "  -5.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.34"
$ws.Range("E18").Value = "  +8.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.47"
$ws.Range("E19").Value = "  +6.89%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "261.03"
$ws.Range("E20").Value = "  +4.67%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0752"
$ws.Range("E21").Value = "  +4.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.52"
$ws.Range("E23").Value = "  +3.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.39"
$ws.Range("E24").Value = "  +2.03%  "

# Row 25
$ws.Range("E25").Value = "  -0.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.08"
$ws.Range("E26").Value = "  +2.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.81"
$ws.Range("E27").Value = "  +5.34%  "

# Row 28
$ws.Range("E28").Value = "  +4.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("E29").Value = "  +5.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.82"
$ws.Range("E31").Value = "  +9.20%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0516"
$ws.Range("E32").Value = "  +3.26%  "

# Row 33
$ws.Range("E33").Value = "  +6.67%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.57"
$ws.Range("E34").Value = "  +8.09%  "

# Row 35
$ws.Range("D35").Value = "1.584.08"
$ws.Range("E35").Value = "  +4.64%  "

# Row 36
$ws.Range("E36").Value = "  +5.06%  "

# Row 37
$ws.Range("E37").Value = "  +2.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "85.44"
$ws.Range("E38").Value = "  +7.30%  "

# Row 39
$ws.Range("E39").Value = "  +7.82%  "

# Row 40
$ws.Range("E40").Value = "  +5.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.84"
$ws.Range("E41").Value = "  +4.32%  "

# Row 42
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.34"
$ws.Range("E42").Value = "  +1.09%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.917"
$ws.Range("E43").Value = "  +7.35%  "

# Row 44
$ws.Range("E44").Value = "  +7.01%  "

# Row 45
$ws.Range("E45").Value = "  +3.51%  "

# Row 46
$ws.Range("E46").Value = "  +4.11%  "

# Row 47
$ws.Range("E47").Value = "  +8.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.48"
$ws.Range("E48").Value = "  +1.81%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.74"
$ws.Range("E49").Value = "  +5.50%  "

# Row 50
$ws.Range("E50").Value = "  +0.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.27"
$ws.Range("E51").Value = "  +16.70%  "
